$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.653.62"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.971.51"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.00"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.41"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.379"
$ws.Range("E9").Value = "  +2.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0789"
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.30"
$ws.Range("E12").Value = "  +4.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.849"
$ws.Range("E13").Value = "  +3.13%  "
$ws.Range("D14").Value = "2.256.11"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.72"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.31"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "1.961.03"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "36.578.07"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.62"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "0.0₃0855"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.11"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.59"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.145"
$ws.Range("E26").Value = "  +6.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.18"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.53"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.38"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +18.00%  "
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.83"
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0616"
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("E34").Value = "  +7.16%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.29"
$ws.Range("E35").Value = "  +3.18%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.38"
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.44"
$ws.Range("E39").Value = "  -13.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0971"
$ws.Range("E40").Value = "  -3.30%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.17"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0210"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.01"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("D45").Value = "1.373.54"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.05"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.83"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.28"
$ws.Range("E50").Value = "  +6.00%  "
$ws.Range("D51").Value = "2.149.98"
$ws.Range("E51").Value = "  +0.74%  "
